$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price (column D) value, kept as text to preserve formatting
$priceUpdates = @{
    2  = "236.90"
    3  = "21.89"
    4  = "5.349"
    6  = "6.479"
    7  = "3.352"
    9  = "1.044"
    10 = "0.1388"
    11 = "0.07332"
    12 = "0.03158"
    13 = "0.02966"
    14 = "0.09244"
    15 = "0.001681"
    16 = "3.258"
    19 = "0.006233"
    20 = "0.005058"
    21 = "0.001051"
    24 = "3.950"
    40 = "0.04106"
    41 = "0.007033"
    44 = "0.008801"
    45 = "0.00005438"
    47 = "0.6760"
    48 = "0.03521"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Map of row -> new Volume(1h) value (column E) - plain text, no leading digits
# that Excel would mis-parse as a number, so no special NumberFormat needed.
$volumeUpdates = @{
    20 = "19HotbitTokenHTBBestin24h"
    24 = "23LEOLEO"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
